# Add two new rows to the "address" sheet duplicating existing rows 2 and 3,
# then make "address" the active sheet with a multi-cell selection.

$wb = $excel.ActiveWorkbook

$customer = $wb.Worksheets.Item("customer")
$address  = $wb.Worksheets.Item("address")

# Append row 4 (duplicate of row 2: eon / kharadi / pune / maharashtra / india)
$address.Range("A4").Value = 1
$address.Range("B4").Value = "eon"
$address.Range("C4").Value = "kharadi"
$address.Range("D4").Value = "pune"
$address.Range("E4").Value = "maharashtra"
$address.Range("F4").Value = "india"

# Append row 5 (duplicate of row 3: eon / magarpatta / pune / maharashtra / india)
$address.Range("A5").Value = 1
$address.Range("B5").Value = "eon"
$address.Range("C5").Value = "magarpatta"
$address.Range("D5").Value = "pune"
$address.Range("E5").Value = "maharashtra"
$address.Range("F5").Value = "india"

# Reset the selection on "customer" sheet
$customer.Range("A2:F2").Select()

# Make "address" the active sheet with its own selection
$address.Activate()
$address.Range("A4:F5").Select()
